$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores values that look numeric (e.g. "304.77" or
# "23.415.35") but are actually plain text in the source workbook. Assigning a
# plain numeric-looking string straight to .Value lets Excel auto-convert it to
# a real number (losing trailing zeros, turning it into scientific notation,
# etc.), so we temporarily force the cell to Text format, assign the string,
# then restore the original (default) style so no stray formatting is left
# behind. The "Volume(1h)" column (E) already contains "%" and spaces, so Excel
# keeps those as text without any special handling.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.415.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.639.85'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3733'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '52.31'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.250'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08117'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("E12").Value = '  -0.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.82'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.590'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001270'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.283'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.632.28'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06895'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.40%  '
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.502'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.41%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '23.424.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.86%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.76'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.067'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.414'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.53'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.334'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '135.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.71%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.280'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.810.47'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.54%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.793'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9512'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.95%  '
$ws.Range("E35").Value = '  +2.38%  '
$ws.Range("E36").Value = '  +0.61%  '
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.07219'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.71%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.091'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.08750'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.79%  '
$ws.Range("E41").Value = '  -1.64%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7044'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6513'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.51%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.329'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07965'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.37'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.24%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.198'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
